$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "Datos actualizados a 6 de Julio de 2020 a las 13:16"  # A1: 'Datos actualizados a 6 de Julio de 2020 a las 11:59' -> 'Datos actualizados a 6 de Julio de 2020 a las 13:16'
$ws.Cells.Item(6,2).Value = 700724  # B6: 699402 -> 700724
$ws.Cells.Item(6,3).Value = 2888  # C6: 1566 -> 2888
$ws.Cells.Item(6,4).Value = 425568  # D6: 425070 -> 425568
$ws.Cells.Item(6,5).Value = 255442  # E6: 254625 -> 255442
$ws.Cells.Item(6,7).Value = 14  # G6: 7 -> 14
$ws.Cells.Item(6,8).Value = 19714  # H6: 19707 -> 19714
$ws.Cells.Item(13,1).Value = "Iran"  # A13: 'Italia' -> 'Iran'
$ws.Cells.Item(13,2).Value = 243051  # B13: 241611 -> 243051
$ws.Cells.Item(13,3).Value = 2613  # C13: 0 -> 2613
$ws.Cells.Item(13,4).Value = 204083  # D13: 192108 -> 204083
$ws.Cells.Item(13,5).Value = 27237  # E13: 14642 -> 27237
$ws.Cells.Item(13,7).Value = 160  # G13: 0 -> 160
$ws.Cells.Item(13,8).Value = 11731  # H13: 34861 -> 11731
$ws.Cells.Item(14,1).Value = "Italia"  # A14: 'Iran' -> 'Italia'
$ws.Cells.Item(14,2).Value = 241611  # B14: 240438 -> 241611
$ws.Cells.Item(14,4).Value = 192108  # D14: 201330 -> 192108
$ws.Cells.Item(14,5).Value = 14642  # E14: 27537 -> 14642
$ws.Cells.Item(14,8).Value = 34861  # H14: 11571 -> 34861
$ws.Cells.Item(39,1).Value = "Oman"  # A39: 'Filipinas' -> 'Oman'
$ws.Cells.Item(39,2).Value = 47735  # B39: 46333 -> 47735
$ws.Cells.Item(39,3).Value = 1557  # C39: 2079 -> 1557
$ws.Cells.Item(39,4).Value = 29146  # D39: 12185 -> 29146
$ws.Cells.Item(39,5).Value = 18371  # E39: 32845 -> 18371
$ws.Cells.Item(39,7).Value = 5  # G39: 6 -> 5
$ws.Cells.Item(39,8).Value = 218  # H39: 1303 -> 218
$ws.Cells.Item(40,1).Value = "Filipinas"  # A40: 'Oman' -> 'Filipinas'
$ws.Cells.Item(40,2).Value = 46333  # B40: 46178 -> 46333
$ws.Cells.Item(40,3).Value = 2079  # C40: 0 -> 2079
$ws.Cells.Item(40,4).Value = 12185  # D40: 27917 -> 12185
$ws.Cells.Item(40,5).Value = 32845  # E40: 18048 -> 32845
$ws.Cells.Item(40,7).Value = 6  # G40: 0 -> 6
$ws.Cells.Item(40,8).Value = 1303  # H40: 213 -> 1303
$ws.Cells.Item(48,2).Value = 32315  # B48: 32268 -> 32315
$ws.Cells.Item(48,3).Value = 47  # C48: 0 -> 47
$ws.Cells.Item(48,5).Value = 1050  # E48: 1003 -> 1050
$ws.Cells.Item(50,5).Value = 4620  # E50: 4621 -> 4620
$ws.Cells.Item(50,7).Value = 1  # G50: 0 -> 1
$ws.Cells.Item(50,8).Value = 98  # H50: 97 -> 98
$ws.Cells.Item(51,2).Value = 29223  # B51: 28973 -> 29223
$ws.Cells.Item(51,3).Value = 250  # C51: 0 -> 250
$ws.Cells.Item(51,4).Value = 20213  # D51: 20026 -> 20213
$ws.Cells.Item(51,5).Value = 7242  # E51: 7197 -> 7242
$ws.Cells.Item(51,7).Value = 18  # G51: 0 -> 18
$ws.Cells.Item(51,8).Value = 1768  # H51: 1750 -> 1768
$ws.Cells.Item(60,2).Value = 18365  # B60: 18280 -> 18365
$ws.Cells.Item(60,3).Value = 85  # C60: 0 -> 85
$ws.Cells.Item(60,4).Value = 16647  # D60: 16615 -> 16647
$ws.Cells.Item(60,5).Value = 1012  # E60: 959 -> 1012
$ws.Cells.Item(63,1).Value = "Nepal"  # A63: 'Argelia' -> 'Nepal'
$ws.Cells.Item(63,2).Value = 15964  # B63: 15941 -> 15964
$ws.Cells.Item(63,3).Value = 180  # C63: 0 -> 180
$ws.Cells.Item(63,4).Value = 6811  # D63: 11492 -> 6811
$ws.Cells.Item(63,5).Value = 9118  # E63: 3497 -> 9118
$ws.Cells.Item(63,7).Value = 1  # G63: 0 -> 1
$ws.Cells.Item(63,8).Value = 35  # H63: 952 -> 35
$ws.Cells.Item(64,1).Value = "Argelia"  # A64: 'Nepal' -> 'Argelia'
$ws.Cells.Item(64,2).Value = 15941  # B64: 15784 -> 15941
$ws.Cells.Item(64,4).Value = 11492  # D64: 6547 -> 11492
$ws.Cells.Item(64,5).Value = 3497  # E64: 9203 -> 3497
$ws.Cells.Item(64,8).Value = 952  # H64: 34 -> 952
$ws.Cells.Item(76,1).Value = "El Salvador"  # A76: 'Kenia' -> 'El Salvador'
$ws.Cells.Item(76,2).Value = 8027  # B76: 7886 -> 8027
$ws.Cells.Item(76,3).Value = 250  # C76: 0 -> 250
$ws.Cells.Item(76,4).Value = 4730  # D76: 2287 -> 4730
$ws.Cells.Item(76,5).Value = 3074  # E76: 5439 -> 3074
$ws.Cells.Item(76,7).Value = 6  # G76: 0 -> 6
$ws.Cells.Item(76,8).Value = 223  # H76: 160 -> 223
$ws.Cells.Item(77,1).Value = "Kenia"  # A77: 'El Salvador' -> 'Kenia'
$ws.Cells.Item(77,2).Value = 7886  # B77: 7777 -> 7886
$ws.Cells.Item(77,4).Value = 2287  # D77: 4588 -> 2287
$ws.Cells.Item(77,5).Value = 5439  # E77: 2966 -> 5439
$ws.Cells.Item(77,7).Value = 0  # G77: 6 -> 0
$ws.Cells.Item(77,8).Value = 160  # H77: 223 -> 160
$ws.Cells.Item(79,2).Value = 7432  # B79: 7411 -> 7432
$ws.Cells.Item(79,3).Value = 21  # C79: 0 -> 21
$ws.Cells.Item(79,4).Value = 3226  # D79: 3184 -> 3226
$ws.Cells.Item(79,5).Value = 4024  # E79: 4045 -> 4024
$ws.Cells.Item(100,1).Value = "Madagascar"  # A100: 'Tailandia' -> 'Madagascar'
$ws.Cells.Item(100,2).Value = 3250  # B100: 3195 -> 3250
$ws.Cells.Item(100,3).Value = 309  # C100: 5 -> 309
$ws.Cells.Item(100,4).Value = 1135  # D100: 3072 -> 1135
$ws.Cells.Item(100,5).Value = 2082  # E100: 65 -> 2082
$ws.Cells.Item(100,7).Value = 1  # G100: 0 -> 1
$ws.Cells.Item(100,8).Value = 33  # H100: 58 -> 33
$ws.Cells.Item(101,1).Value = "Tailandia"  # A101: 'Croacia' -> 'Tailandia'
$ws.Cells.Item(101,2).Value = 3195  # B101: 3151 -> 3195
$ws.Cells.Item(101,3).Value = 5  # C101: 0 -> 5
$ws.Cells.Item(101,4).Value = 3072  # D101: 2196 -> 3072
$ws.Cells.Item(101,5).Value = 65  # E101: 842 -> 65
$ws.Cells.Item(101,8).Value = 58  # H101: 113 -> 58
$ws.Cells.Item(102,1).Value = "Croacia"  # A102: 'Guinea Ecuatorial' -> 'Croacia'
$ws.Cells.Item(102,2).Value = 3151  # B102: 3071 -> 3151
$ws.Cells.Item(102,4).Value = 2196  # D102: 842 -> 2196
$ws.Cells.Item(102,5).Value = 842  # E102: 2178 -> 842
$ws.Cells.Item(102,8).Value = 113  # H102: 51 -> 113
$ws.Cells.Item(103,1).Value = "Guinea Ecuatorial"  # A103: 'Somalia' -> 'Guinea Ecuatorial'
$ws.Cells.Item(103,2).Value = 3071  # B103: 2997 -> 3071
$ws.Cells.Item(103,4).Value = 842  # D103: 1014 -> 842
$ws.Cells.Item(103,5).Value = 2178  # E103: 1891 -> 2178
$ws.Cells.Item(103,8).Value = 51  # H103: 92 -> 51
$ws.Cells.Item(104,1).Value = "Somalia"  # A104: 'Albania' -> 'Somalia'
$ws.Cells.Item(104,2).Value = 2997  # B104: 2964 -> 2997
$ws.Cells.Item(104,3).Value = 0  # C104: 71 -> 0
$ws.Cells.Item(104,4).Value = 1014  # D104: 1702 -> 1014
$ws.Cells.Item(104,5).Value = 1891  # E104: 1183 -> 1891
$ws.Cells.Item(104,7).Value = 0  # G104: 3 -> 0
$ws.Cells.Item(104,8).Value = 92  # H104: 79 -> 92
$ws.Cells.Item(105,1).Value = "Albania"  # A105: 'Madagascar' -> 'Albania'
$ws.Cells.Item(105,2).Value = 2964  # B105: 2941 -> 2964
$ws.Cells.Item(105,3).Value = 71  # C105: 0 -> 71
$ws.Cells.Item(105,4).Value = 1702  # D105: 1108 -> 1702
$ws.Cells.Item(105,5).Value = 1183  # E105: 1801 -> 1183
$ws.Cells.Item(105,7).Value = 3  # G105: 0 -> 3
$ws.Cells.Item(105,8).Value = 79  # H105: 32 -> 79
$ws.Cells.Item(112,4).Value = 1917  # D112: 1903 -> 1917
$ws.Cells.Item(112,5).Value = 148  # E112: 162 -> 148
$ws.Cells.Item(137,1).Value = "Burkina Faso"  # A137: 'Suazilandia' -> 'Burkina Faso'
$ws.Cells.Item(137,2).Value = 1000  # B137: 988 -> 1000
$ws.Cells.Item(137,3).Value = 13  # C137: 0 -> 13
$ws.Cells.Item(137,4).Value = 858  # D137: 547 -> 858
$ws.Cells.Item(137,5).Value = 89  # E137: 428 -> 89
$ws.Cells.Item(137,8).Value = 53  # H137: 13 -> 53
$ws.Cells.Item(138,1).Value = "Suazilandia"  # A138: 'Mozambique' -> 'Suazilandia'
$ws.Cells.Item(138,2).Value = 988  # B138: 987 -> 988
$ws.Cells.Item(138,4).Value = 547  # D138: 256 -> 547
$ws.Cells.Item(138,5).Value = 428  # E138: 723 -> 428
$ws.Cells.Item(138,8).Value = 13  # H138: 8 -> 13
$ws.Cells.Item(139,1).Value = "Mozambique"  # A139: 'Burkina Faso' -> 'Mozambique'
$ws.Cells.Item(139,4).Value = 256  # D139: 854 -> 256
$ws.Cells.Item(139,5).Value = 723  # E139: 80 -> 723
$ws.Cells.Item(139,8).Value = 8  # H139: 53 -> 8
$ws.Cells.Item(160,2).Value = 369  # B160: 355 -> 369
$ws.Cells.Item(160,3).Value = 14  # C160: 0 -> 14
$ws.Cells.Item(160,5).Value = 29  # E160: 15 -> 29
